$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("G34").Value = 1.3
$ws.Range("H34").Value = 5
$ws.Range("I34").Value = 11
$ws.Range("J34").Value = 1.73
$ws.Range("K34").Value = 2.6
$ws.Range("L34").Value = 8.5
$ws.Range("M34").Value = 1.04
$ws.Range("N34").Value = 13
$ws.Range("O34").Value = 1.2
$ws.Range("P34").Value = 4.33
$ws.Range("Q34").Value = 1.65
$ws.Range("R34").Value = 2.2
$ws.Range("S34").Value = 1.3
$ws.Range("T34").Value = 3.4
$ws.Range("U34").Value = 2.1
$ws.Range("V34").Value = 1.67
$ws.Range("W34").Value = 7.5
$ws.Range("Z34").Value = 8
$ws.Range("AC34").Value = 13
$ws.Range("AD34").Value = 9.5
$ws.Range("AE34").Value = 21
$ws.Range("AG34").Value = 301
$ws.Range("AH34").Value = 23
$ws.Range("AJ34").Value = 29
$ws.Range("AK34").Value = 126
$ws.Range("AL34").Value = 67
$ws.Range("AO34").Value = 6
$ws.Range("AT34").Value = 3.4
$ws.Range("AU34").Value = 9.5
$ws.Range("AV34").Value = 67

# Row 43
$ws.Range("G43").Value = 2.25
$ws.Range("H43").Value = 3.3
$ws.Range("I43").Value = 3.2
$ws.Range("J43").Value = 2.88
$ws.Range("L43").Value = 3.5
$ws.Range("Q43").Value = 1.8
$ws.Range("R43").Value = 2
$ws.Range("S43").Value = 1.36
$ws.Range("T43").Value = 3
$ws.Range("W43").Value = 9
$ws.Range("Y43").Value = 9.5
$ws.Range("AN43").Value = 4.5
$ws.Range("AP43").Value = 21
$ws.Range("AT43").Value = 3
$ws.Range("AX43").Value = 5

# Row 49
$ws.Range("G49").Value = 2
$ws.Range("I49").Value = 4.2
$ws.Range("L49").Value = 4.75
$ws.Range("Q49").Value = 2.6
$ws.Range("R49").Value = 1.48
$ws.Range("Z49").Value = 17
$ws.Range("AA49").Value = 21
$ws.Range("AE49").Value = 19
$ws.Range("AH49").Value = 9
$ws.Range("AI49").Value = 19
$ws.Range("AJ49").Value = 15
$ws.Range("AK49").Value = 41
$ws.Range("AO49").Value = 12
$ws.Range("AX49").Value = 5.5
$ws.Range("AY49").Value = 23
$ws.Range("BA49").Value = 81
$ws.Range("BB49").Value = 126

# Row 50
$ws.Range("G50").Value = 1.53
$ws.Range("H50").Value = 3.7
$ws.Range("I50").Value = 7
$ws.Range("J50").Value = 2.2
$ws.Range("K50").Value = 2.05
$ws.Range("L50").Value = 7.5
$ws.Range("M50").Value = 1.11
$ws.Range("N50").Value = 6.5
$ws.Range("U50").Value = 2.63
$ws.Range("V50").Value = 1.44
$ws.Range("Z50").Value = 10
$ws.Range("AD50").Value = 7.5
$ws.Range("AH50").Value = 12
$ws.Range("AI50").Value = 34
$ws.Range("AJ50").Value = 23
$ws.Range("AL50").Value = 67
$ws.Range("AS50").Value = 301
$ws.Range("AV50").Value = 101
$ws.Range("AX50").Value = 8
$ws.Range("AZ50").Value = 51
$ws.Range("BA50").Value = 201

# Row 62
$ws.Range("G62").Value = 1.36
$ws.Range("H62").Value = 4.5
$ws.Range("I62").Value = 9.5
$ws.Range("L62").Value = 8
$ws.Range("Q62").Value = 1.83
$ws.Range("R62").Value = 1.98
$ws.Range("Z62").Value = 8.5
$ws.Range("AC62").Value = 10
$ws.Range("AD62").Value = 8.5
$ws.Range("AE62").Value = 21
$ws.Range("AF62").Value = 67
$ws.Range("AH62").Value = 21
$ws.Range("AP62").Value = 21
$ws.Range("AQ62").Value = 19
$ws.Range("AX62").Value = 9

# Row 97
$ws.Range("G97").Value = 2.1
$ws.Range("H97").Value = 3.3
$ws.Range("I97").Value = 3.7
$ws.Range("J97").Value = 2.75
$ws.Range("K97").Value = 2.1
$ws.Range("L97").Value = 4
$ws.Range("O97").Value = 1.3
$ws.Range("P97").Value = 3.4
$ws.Range("Q97").Value = 2.05
$ws.Range("R97").Value = 1.8
$ws.Range("U97").Value = 1.75
$ws.Range("V97").Value = 2
$ws.Range("X97").Value = 10
$ws.Range("Y97").Value = 9
$ws.Range("Z97").Value = 19
$ws.Range("AA97").Value = 17
$ws.Range("AC97").Value = 9.5
$ws.Range("AD97").Value = 6
$ws.Range("AE97").Value = 13
$ws.Range("AG97").Value = 201
$ws.Range("AH97").Value = 11
$ws.Range("AI97").Value = 19
$ws.Range("AL97").Value = 29
$ws.Range("AM97").Value = 34
$ws.Range("AO97").Value = 11
$ws.Range("AQ97").Value = 41
$ws.Range("AY97").Value = 19
$ws.Range("BB97").Value = 81

# Row 109
$ws.Range("M109").Value = 1.07
$ws.Range("N109").Value = 8.5
$ws.Range("Z109").Value = 9.5

# Row 173
$ws.Range("O173").Value = 1.22
$ws.Range("P173").Value = 4

# Row 175
$ws.Range("G175").Value = 1.93
$ws.Range("H175").Value = 3.3
$ws.Range("I175").Value = 3.55
$ws.Range("J175").Value = 2.55
$ws.Range("K175").Value = 2.12
$ws.Range("L175").Value = 4.1
$ws.Range("M175").Value = 1.07
$ws.Range("N175").Value = 6.8
$ws.Range("O175").Value = 1.35
$ws.Range("P175").Value = 2.95
$ws.Range("Q175").Value = 2.05
$ws.Range("R175").Value = 1.72
$ws.Range("T175").Value = 2.67
$ws.Range("U175").Value = 1.88
$ws.Range("W175").Value = 6.7
$ws.Range("X175").Value = 8.75
$ws.Range("Y175").Value = 8.5
$ws.Range("Z175").Value = 16.5
$ws.Range("AA175").Value = 16.5
$ws.Range("AB175").Value = 30
$ws.Range("AC175").Value = 6.8
$ws.Range("AD175").Value = 6.5
$ws.Range("AE175").Value = 16
$ws.Range("AH175").Value = 9.5
$ws.Range("AI175").Value = 18.5
$ws.Range("AJ175").Value = 13
$ws.Range("AK175").Value = 50
$ws.Range("AL175").Value = 35
$ws.Range("AM175").Value = 45
$ws.Range("AN175").Value = 3.8
$ws.Range("AO175").Value = 10
$ws.Range("AP175").Value = 19.5
$ws.Range("AQ175").Value = 37
$ws.Range("AR175").Value = 75
$ws.Range("AS175").Value = 250
$ws.Range("AT175").Value = 2.67
$ws.Range("AU175").Value = 7.5
$ws.Range("AX175").Value = 5.4
$ws.Range("AY175").Value = 20
$ws.Range("AZ175").Value = 29
$ws.Range("BA175").Value = 110
$ws.Range("BB175").Value = 150
$ws.Range("BC175").Value = 400

# Row 177
$ws.Range("G177").Value = 3.2
$ws.Range("H177").Value = 2.85
$ws.Range("I177").Value = 2.32
$ws.Range("J177").Value = 3.9
$ws.Range("K177").Value = 1.85
$ws.Range("L177").Value = 3.1
$ws.Range("M177").Value = 1.1
$ws.Range("N177").Value = 7.02
$ws.Range("S177").Value = 1.55
$ws.Range("T177").Value = 2.15
$ws.Range("W177").Value = 7.7
$ws.Range("X177").Value = 15.5
$ws.Range("Y177").Value = 11.75
$ws.Range("Z177").Value = 45
$ws.Range("AA177").Value = 35
$ws.Range("AB177").Value = 50
$ws.Range("AC177").Value = 6.4
$ws.Range("AD177").Value = 5.7
$ws.Range("AH177").Value = 5.9
$ws.Range("AJ177").Value = 9.75
$ws.Range("AK177").Value = 24
$ws.Range("AL177").Value = 24
$ws.Range("AM177").Value = 45
$ws.Range("AN177").Value = 4.85
$ws.Range("AO177").Value = 19
$ws.Range("AP177").Value = 30
$ws.Range("AQ177").Value = 110
$ws.Range("AR177").Value = 175
$ws.Range("AT177").Value = 2.1
$ws.Range("AU177").Value = 7.8
$ws.Range("AV177").Value = 90
$ws.Range("AX177").Value = 4
$ws.Range("AY177").Value = 13.5
$ws.Range("AZ177").Value = 27
$ws.Range("BA177").Value = 65
$ws.Range("BB177").Value = 120
$ws.Range("BC177").Value = 450

